# Adds the "Assign Alfresco Folder" Drools rule row to the Costsheet
# rules sheet (Sheet1), matching the author's commit:
#   "added drools rules added costsheet xsl"
#
# Functional changes reproduced:
#   - New rule row 23: Rule Name / When Expression is True / Set Field Value
#   - Column D widened to fit the longer action expression
#   - Selection/active cell moved to B23 (next empty rule row)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New rule row (row 23) ------------------------------------------------
$ws.Range("B23").Value = "Assign Alfresco Folder"
$ws.Range("C23").Value = "container?.folder?.cmisFolderId == null"
$ws.Range("D23").Value = "setEcmFolderPath, '/Sites/acm/documentLibrary/Expenses/' + dateFormat('yyyyMMdd') + '_' + `$acmCostsheet.getId()"

# --- Column D width: 116 chars wide (was 100.42578125) --------------------
# ColumnWidth setter adds Excel's standard 5/6-character padding, so request
# 116 - 5/6 to land exactly on a stored width of 116.
$ws.Columns.Item(4).ColumnWidth = 115.16666666666667

# --- Selection moves to the newly-added rule row ---------------------------
$ws.Range("B23").Select()
